$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the "İsim" column to B.
$ws.Range("A1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("A1").Value = "Numara"

# Student numbers for the newly inserted column A.
$ws.Range("A2").Value = 201012
$ws.Range("A3").Value = 201013
$ws.Range("A4").Value = 201014
$ws.Range("A5").Value = 201015

# Reset the active selection to A5, matching the saved workbook view.
$ws.Range("A5").Select()
